$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"

$ws.Range("C3").Value = "Desenho Técnico - MCT-1A"
$ws.Range("C4").Value = "Desenho Técnico - MCT-1A"
$ws.Range("C6").Value = "Desenho Técnico - ELT-1A"
